# Add line data to all facility line responses.
#
# - Refresh the uuid column (G) on the existing rows (2,3) to the new uuid.
# - Append two new rows (4,5) that replicate rows 2 and 3 (same line/date/
#   descr/start/end/gap/downtime) but stamped with the new uuid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUuid = "f72d2039-1438-4515-aa58-11f1078e2401"

# --- Existing rows: just refresh the uuid column ---
$ws.Range("G2").Value = $newUuid
$ws.Range("G3").Value = $newUuid

# --- New rows 4 (mirrors row 2) and 5 (mirrors row 3) ---
$sourceRows = @(2, 3)
$targetRows = @(4, 5)

for ($i = 0; $i -lt $sourceRows.Length; $i++) {
    $src = $sourceRows[$i]
    $dst = $targetRows[$i]

    # A: plain text, no special format on the source cell.
    $ws.Range("A$dst").Value = $ws.Range("A$src").Value()

    # B: looks like a date ("01/09/2024") - force it to stay text (matching
    # the source, which has no numeric/date style) using a leading quote,
    # then strip the quote-prefix style it leaves behind so the cell ends
    # up back on the default (unstyled) format, same as the source cell.
    $ws.Range("B$dst").Value = "'" + $ws.Range("B$src").Value()
    $ws.Range("B$dst").Style = "Normal"

    # C: plain text, no special format on the source cell.
    $ws.Range("C$dst").Value = $ws.Range("C$src").Value()

    # D, E: numeric date/time values - copy the source's number format so
    # the new cells land on the same style as the source (s="2").
    $ws.Range("D$dst").NumberFormat = $ws.Range("D$src").NumberFormat
    $ws.Range("D$dst").Value = $ws.Range("D$src").Value()

    $ws.Range("E$dst").NumberFormat = $ws.Range("E$src").NumberFormat
    $ws.Range("E$dst").Value = $ws.Range("E$src").Value()

    # F: plain numeric, no special format on the source cell.
    $ws.Range("F$dst").Value = $ws.Range("F$src").Value()

    # G: uuid column - new rows get the new uuid.
    $ws.Range("G$dst").Value = $newUuid

    # H: downtime - copy the source's number format (s="3").
    $ws.Range("H$dst").NumberFormat = $ws.Range("H$src").NumberFormat
    $ws.Range("H$dst").Value = $ws.Range("H$src").Value()
}
